# Updated Indonesia files compatible with v3.3.1
# Applies the content edits to the "EoDSDwSP" elasticity workbook:
#  - About sheet: fill in the previously-empty note cell (A9) with "Notes:"
#    and add two new note lines below it (A10, A11).
#  - EoDSDwSP sheet: clarify the "Elasticity" label to "Elasticity (dimensionless)".

$wb = $excel.ActiveWorkbook

# --- "About" sheet: add Notes section ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A9").Value = "Notes:"
$wsAbout.Range("A9").Style = "Normal"
$wsAbout.Cells.Item(9, 1).Font.Bold = $true

$wsAbout.Range("A10").Value = "Elasticities intended to reflect change in deployment with changing"
$wsAbout.Range("A11").Value = "distributed solar price (through subsidies)."

# --- "EoDSDwSP" sheet: clarify Elasticity label ---
$wsResult = $wb.Worksheets.Item("EoDSDwSP")
$wsResult.Range("B1").Value = "Elasticity (dimensionless)"
[void]$wsResult.Range("B2").Select()

# Leave the "About" sheet active (matches the source file, which keeps
# tabSelected on "About") with its own remembered selection at A12.
[void]$wsAbout.Activate()
[void]$wsAbout.Range("A12").Select()
